# Updated cryptos list on Thu Feb  1 16:51:06 UTC 2024 with GitHub Actions
# Re-applies the latest coinranking.com price/volume snapshot onto Sheet1 (A1:E51).
# Cells that would otherwise auto-parse as a pure number (e.g. "1.73") are written
# with a leading apostrophe so they stay text, matching the sheet's original inlineStr layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.590.20"
$ws.Range("E2").Value = "  -2.03%  "
# Row 3
$ws.Range("D3").Value = "2.287.53"
$ws.Range("E3").Value = "  -2.22%  "
# Row 4
$ws.Range("D4").Value = "'0.999"
# Row 5
$ws.Range("D5").Value = "'299.79"
$ws.Range("E5").Value = "  -1.75%  "
# Row 6
$ws.Range("D6").Value = "'96.18"
$ws.Range("E6").Value = "  -5.13%  "
# Row 7
$ws.Range("D7").Value = "'0.500"
$ws.Range("E7").Value = "  -1.99%  "
# Row 8
$ws.Range("E8").Value = "  +0.01%  "
# Row 9
$ws.Range("D9").Value = "'0.491"
$ws.Range("E9").Value = "  -3.29%  "
# Row 10
$ws.Range("D10").Value = "'33.26"
$ws.Range("E10").Value = "  -5.97%  "
# Row 11
$ws.Range("D11").Value = "'0.0791"
$ws.Range("E11").Value = "  -0.99%  "
# Row 12
$ws.Range("D12").Value = "'49.05"
$ws.Range("E12").Value = "  -4.07%  "
# Row 13
$ws.Range("E13").Value = "  +1.99%  "
# Row 14
$ws.Range("D14").Value = "'16.80"
$ws.Range("E14").Value = "  +6.87%  "
# Row 15
$ws.Range("D15").Value = "'6.75"
$ws.Range("E15").Value = "  -0.96%  "
# Row 16
$ws.Range("D16").Value = "2.637.35"
$ws.Range("E16").Value = "  -1.94%  "
# Row 17
$ws.Range("D17").Value = "2.257.01"
$ws.Range("E17").Value = "  -2.80%  "
# Row 18
$ws.Range("D18").Value = "'0.794"
$ws.Range("E18").Value = "  -1.83%  "
# Row 19
$ws.Range("D19").Value = "42.433.08"
$ws.Range("E19").Value = "  -2.08%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0895"
$ws.Range("E20").Value = "  -1.73%  "
# Row 21
$ws.Range("D21").Value = "'11.47"
$ws.Range("E21").Value = "  -3.66%  "
# Row 22
$ws.Range("D22").Value = "'5.98"
$ws.Range("E22").Value = "  -2.26%  "
# Row 23
$ws.Range("D23").Value = "'66.78"
$ws.Range("E23").Value = "  -2.28%  "
# Row 24
$ws.Range("D24").Value = "'235.88"
$ws.Range("E24").Value = "  -0.89%  "
# Row 25
$ws.Range("D25").Value = "'1.96"
$ws.Range("E25").Value = "  -1.55%  "
# Row 26
$ws.Range("E26").Value = "  +0.08%  "
# Row 27
$ws.Range("D27").Value = "'2.45"
$ws.Range("E27").Value = "  -3.31%  "
# Row 28
$ws.Range("D28").Value = "'24.29"
$ws.Range("E28").Value = "  -3.16%  "
# Row 29
$ws.Range("E29").Value = "  -0.89%  "
# Row 30
$ws.Range("D30").Value = "'166.55"
$ws.Range("E30").Value = "  +0.43%  "
# Row 31
$ws.Range("D31").Value = "'33.65"
$ws.Range("E31").Value = "  -3.04%  "
# Row 32
$ws.Range("D32").Value = "'9.09"
$ws.Range("E32").Value = "  -1.72%  "
# Row 33
$ws.Range("D33").Value = "'0.999"
$ws.Range("E33").Value = "  +0.07%  "
# Row 34
$ws.Range("D34").Value = "'4.74"
$ws.Range("E34").Value = "  +3.84%  "
# Row 35
$ws.Range("D35").Value = "'4.94"
$ws.Range("E35").Value = "  -2.78%  "
# Row 36
$ws.Range("D36").Value = "'2.40"
$ws.Range("E36").Value = "  -0.57%  "
# Row 37
$ws.Range("D37").Value = "'16.76"
$ws.Range("E37").Value = "  -0.96%  "
# Row 38
$ws.Range("D38").Value = "'0.0692"
$ws.Range("E38").Value = "  -2.16%  "
# Row 39
$ws.Range("D39").Value = "'2.80"
$ws.Range("E39").Value = "  -4.20%  "
# Row 40
$ws.Range("D40").Value = "'0.0995"
$ws.Range("E40").Value = "  -3.18%  "
# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'1.73"
$ws.Range("E41").Value = "  -5.24%  "
# Row 42
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.109"
$ws.Range("E42").Value = "  -2.11%  "
# Row 43
$ws.Range("D43").Value = "'2.37"
$ws.Range("E43").Value = "  -1.52%  "
# Row 44
$ws.Range("D44").Value = "1.955.60"
$ws.Range("E44").Value = "  -1.17%  "
# Row 45
$ws.Range("D45").Value = "'0.0279"
$ws.Range("E45").Value = "  -1.88%  "
# Row 46
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'17.46"
$ws.Range("E46").Value = "  -6.04%  "
# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").Value = "'9.69"
$ws.Range("E47").Value = "  -3.23%  "
# Row 48
$ws.Range("D48").Value = "'2.81"
$ws.Range("E48").Value = "  -4.41%  "
# Row 49
$ws.Range("D49").Value = "2.505.91"
$ws.Range("E49").Value = "  -1.93%  "
# Row 50
$ws.Range("D50").Value = "'52.46"
$ws.Range("E50").Value = "  -6.67%  "
# Row 51
$ws.Range("D51").Value = "'2.72"
$ws.Range("E51").Value = "  -2.39%  "
